# Regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 0
    3  = 1
    4  = 1
    5  = 2
    6  = 2
    7  = 1
    8  = 2
    9  = 0
    10 = 0
    11 = 1
    12 = 2
    13 = 1
    14 = 3
    15 = 1
    16 = 1
    17 = 0
    18 = 3
    19 = 2
    20 = 1
    21 = 0
    22 = 2
    23 = 1
    24 = 0
    25 = 0
    26 = 2
    27 = 2
    28 = 1
    29 = 1
    30 = 0
    31 = 1
    32 = 0
    33 = 3
    34 = 2
    35 = 2
    36 = 0
    37 = 2
    38 = 2
    39 = 2
    40 = 0
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
